# "(almost) fixed no-username reg"
# Adds two newly-registered users (rows 3 & 4) to the usersInfo sheet and
# mirrors their default equipment loadout (copied from the first user's
# row) into the usersEquipped sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "usersInfo": new user rows -------------------------------------
$wsInfo = $wb.Worksheets.Item("usersInfo")

$wsInfo.Range("A3").Value = "Константин2"
$wsInfo.Range("B3").Value = "2К"
$wsInfo.Range("C3").Value = "26.2.2023"

$wsInfo.Range("A4").Value = "qiviCHAN3"
$wsInfo.Range("B4").Value = "3q"
$wsInfo.Range("C4").Value = "26.2.2023"

# Keep the trailing blank row (row 5) that Excel's save round-trip leaves
# behind after this block of edits - touch a row-only property with a
# value equal to its default so no cell content/formatting is introduced.
$wsInfo.Rows.Item(5).OutlineLevel = 0

# --- Sheet "usersEquipped": default loadout for the new users -------------
$wsEquipped = $wb.Worksheets.Item("usersEquipped")

$weaponCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
$weaponCodes = @("ke0","сс0","sy0","fy0","gt0","sf0","by0","ml0","sr0","as0","se0","je0","bg0","gn0","pm0","vl0","on0","or0")

$wsEquipped.Range("A3").Value = "2К"
for ($i = 0; $i -lt $weaponCols.Length; $i++) {
    $wsEquipped.Range($weaponCols[$i] + "3").Value = $weaponCodes[$i]
}

$wsEquipped.Range("A4").Value = "3q"
for ($i = 0; $i -lt $weaponCols.Length; $i++) {
    $wsEquipped.Range($weaponCols[$i] + "4").Value = $weaponCodes[$i]
}
